# "added the expense overview"
# Updates the Income sheet: row 2 becomes "savings", a new "Youtube" row is
# inserted as row 3, and the old "Salary" row (now row 4) gets new figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 2's formatting (incl. the date-style C column) and insert it as the
# new row 3 -- this shifts the existing row 3 ("Salary") down to row 4 while
# keeping consistent cell styles (s="1" on the date column) for the new row.
$ws.Range("A2:C2").Copy()
$ws.Rows.Item(3).Insert()

# Row 2: "gig" / 8000 / 2025-04-01  ->  "savings" / 4300 / 2025-12-25
$ws.Range("A2").Value = "savings"
$ws.Range("B2").Value = 4300
$ws.Range("C2").Value = 46016.22928240741

# Row 3 (new): "Youtube" / 3000 / 2025-12-23
$ws.Range("A3").Value = "Youtube"
$ws.Range("B3").Value = 3000
$ws.Range("C3").Value = 46014.22928240741

# Row 4 (previously row 3): "Salary" / 4000 / 2025-02-01  ->  "Salary" / 34000 / 2025-12-22
$ws.Range("A4").Value = "Salary"
$ws.Range("B4").Value = 34000
$ws.Range("C4").Value = 46013.22928240741
